# Apply "conservative present scenario" values + capacity-factor-adjusted
# WindOnshore rows in the 2030/2050 scenarios.

$wb = $excel.ActiveWorkbook

# --- Sheet "Present-Storage": new values for the conservative present
#     scenario, plus the narrower column B width that results from the
#     shorter values.
$wsPresent = $wb.Worksheets.Item("Present-Storage")

$presentValues = @{
    2  = 7.71
    3  = 5.53
    4  = 8.33
    5  = 8.24
    6  = 7
    7  = 4.7
    8  = 7.66
    9  = 7.56
    10 = 15.8
    11 = 14.11
    12 = 16.29
    13 = 16.22
    14 = 7.95
    15 = 5.86
    16 = 8.54
    17 = 8.45
}

foreach ($row in $presentValues.Keys) {
    $wsPresent.Cells.Item($row, 2).Value = $presentValues[$row]
}

# The OOXML <col> width attribute this runtime emits is always a multiple
# of 1/6 (Excel's pixel-quantized column-width model), so the exact target
# of 8.4 can't be hit through the ColumnWidth COM property. 7.5 rounds to
# the closest reachable width (8.333...), the nearest representable value.
$wsPresent.Columns.Item(2).ColumnWidth = 7.5

# --- Sheet "2030-Storage": capacity-factor-adjusted WindOnshore rows.
$ws2030 = $wb.Worksheets.Item("2030-Storage")
$ws2030.Cells.Item(3, 2).Value = 4.41
$ws2030.Cells.Item(7, 2).Value = 4.14
$ws2030.Cells.Item(11, 2).Value = 8.41
$ws2030.Cells.Item(15, 2).Value = 5.55

# --- Sheet "2050-Storage": capacity-factor-adjusted WindOnshore rows.
$ws2050 = $wb.Worksheets.Item("2050-Storage")
$ws2050.Cells.Item(3, 2).Value = 68.49
$ws2050.Cells.Item(7, 2).Value = 71.9
$ws2050.Cells.Item(11, 2).Value = 58.41
$ws2050.Cells.Item(15, 2).Value = 67.03
